# "check arguments make sense, right way round etc"
#
# The scenarios table had sigma_true (col E) and mu_cf (col G) values that
# were wrong / implausible for the simulation study. Correct them:
#   - sigma_true of 1  -> 0.4
#   - sigma_true of 10 -> 0.7
#   - mu_cf of 0.2     -> -1.39
#
# Also reflects the reviewer's on-screen state while looking this over:
# zoomed to 90% and with the bottom block of corrected prop_censoring /
# sigma_true values (E26:E33) selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33

$sigmaCol = 5   # E: sigma_true
$muCfCol  = 7   # G: mu_cf

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sigmaCell = $ws.Cells.Item($r, $sigmaCol)
    $sigmaVal = $sigmaCell.Value2

    if ($sigmaVal -eq 1) {
        $sigmaCell.Value = 0.4
    } elseif ($sigmaVal -eq 10) {
        $sigmaCell.Value = 0.7
    }

    $muCfCell = $ws.Cells.Item($r, $muCfCol)
    $muCfVal = $muCfCell.Value2

    if ($muCfVal -eq 0.2) {
        $muCfCell.Value = -1.39
    }
}

# Reviewer state: zoom to 90% and select the last block of edited rows.
$null = $ws.Activate()
$excel.ActiveWindow.Zoom = 90
$null = $ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("E26:E33").Select()
